$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect the new date
$ws.Name = "Through 2022-08-25"

# Update the August row label
$ws.Range("A9").Value = "August (through 08-25)"

# Update August row values (row 9), columns C:I
$ws.Range("C9").Value = 59
$ws.Range("D9").Value = 71
$ws.Range("E9").Value = 46
$ws.Range("F9").Value = 37
$ws.Range("G9").Value = 141
$ws.Range("H9").Value = 129
$ws.Range("I9").Value = 135

# Update Total row values (row 10), columns C:I
$ws.Range("C10").Value = 361
$ws.Range("D10").Value = 536
$ws.Range("E10").Value = 471
$ws.Range("F10").Value = 341
$ws.Range("G10").Value = 762
$ws.Range("H10").Value = 1039
$ws.Range("I10").Value = 1106
